# Atualização automática de preços de eletricidade
# Updates row 2 (the single data row) of the SpotPTTable with the latest
# daily/hourly spot price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial value) - keep the existing date formatting (style already applied)
$ws.Range("A2").Value = 45889

# Hourly prices 0h-1h ... 23h-24h (columns B..Y)
$ws.Range("B2").Value = 95.58
$ws.Range("C2").Value = 88.03
$ws.Range("D2").Value = 83
$ws.Range("E2").Value = 81
$ws.Range("F2").Value = 76.2
$ws.Range("G2").Value = 82
$ws.Range("H2").Value = 88.45999999999999
$ws.Range("I2").Value = 96.54000000000001
$ws.Range("J2").Value = 89.81
$ws.Range("K2").Value = 71.55
$ws.Range("L2").Value = 49.8
$ws.Range("M2").Value = 49.9
$ws.Range("N2").Value = 49.02
$ws.Range("O2").Value = 35
$ws.Range("P2").Value = 31.53
$ws.Range("Q2").Value = 27.2
$ws.Range("R2").Value = 27.99
$ws.Range("S2").Value = 37
$ws.Range("T2").Value = 49.27
$ws.Range("U2").Value = 69.27
$ws.Range("V2").Value = 78.73
$ws.Range("W2").Value = 100.07
$ws.Range("X2").Value = 97.52
$ws.Range("Y2").Value = 84.98

# Daily average price
$ws.Range("Z2").Value = 68.31

# Slot_4h_max ("AA2") is unchanged -> "20h-24h"

# Slot_4h_price
$ws.Range("AB2").Value = 90.31999999999999

# Slot_2h_frist (sic)
$ws.Range("AC2").Value = "6h-8h"

# Slot_2h_frist_price
$ws.Range("AD2").Value = 92.5

# Slot_2h_second
$ws.Range("AE2").Value = "0h-2h"

# Slot_2h_second_price
$ws.Range("AF2").Value = 91.8

# Slot_min_price ("AG2") is unchanged -> "10h-18h"
